# Insert a new weekly price record as row 401, shifting existing rows 401:495
# down to 402:496 (dimension grows from A1:R495 to A1:R496).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(401).Insert()

$ws.Range("A401").Value = 3
$ws.Range("B401").Value = "Femacal de La Calera"
$ws.Range("C401").Value = "Coquimbo"
$ws.Range("D401").Value = 44754
$ws.Range("E401").Value = 5
$ws.Range("F401").Value = 100112032
$ws.Range("G401").Value = "Zapallo italiano"
$ws.Range("H401").Value = "Sin especificar"
$ws.Range("I401").Value = "Primera"
$ws.Range("J401").Value = 278
$ws.Range("K401").Value = 9000
$ws.Range("L401").Value = 10000
$ws.Range("M401").Value = 9514
$ws.Range("N401").Value = "$/caja 70 unidades"
$ws.Range("O401").Value = "Región de Arica y Parinacota"
$ws.Range("P401").Value = 136
$ws.Range("Q401").Value = 70
$ws.Range("R401").Value = "Hortaliza"
